# Generate Report for Handoff
# Update status text from "In Translation" -> "Ready for handoff" and
# bump the "latest" timestamps, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: status columns for zh-cn (E) and de-de (F), plus the
# "Latest HO Xliff Generate Date" column (G).
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-13 06:48:24"

# zh-cn sheet: Status (C) and Latest Handoff Datetime (H)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-13 06:48:16"

# de-de sheet: Status (C) and Latest Handoff Datetime (H)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-13 06:48:24"

# The "Status" column got wider now that the new text is longer than
# "In Translation" - mirror Excel's resulting column width on every sheet
# that shows the status (Overview shows it twice, for zh-cn and de-de).
$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25
$wsZhCn.Columns.Item(3).ColumnWidth = 16.25
$wsDeDe.Columns.Item(3).ColumnWidth = 16.25
